$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Rv1451", 3, "ctaB Rv1451", "FUNCTION: Converts heme B (protoheme IX) to heme O by substitution of the vinyl group on carbon 2 of heme B porphyrin ring with a hydroxyethyl farnesyl side group. {ECO:0000255|HAMAP-Rule:MF_00154}.", 34),
    @("Rv1797", 3, "eccE5 Rv1797", "FUNCTION: Part of the ESX-5 specialized secretion system, which is responsible for the secretion of EsxN and a number of PE_PGRS and PPE proteins, including PPE41. {ECO:0000269|PubMed:22340629, ECO:0000269|PubMed:22925462}.", 34),
    @("Rv1795", 3, "eccD5 Rv1795", "FUNCTION: Part of the ESX-5 specialized secretion system, which is responsible for the secretion of EsxN and a number of PE_PGRS and PPE proteins, including PPE41. {ECO:0000269|PubMed:22340629, ECO:0000269|PubMed:22925462}.", 34),
    @("Rv0290", 3, "eccD3 Rv0290", "FUNCTION: Part of the ESX-3 specialized secretion system, which is important for iron and zinc uptake or homeostasis. {ECO:0000269|PubMed:19684129, ECO:0000269|PubMed:24155985}.", 34),
    @("Rv0987", 3, "Rv0987", "FUNCTION: Probably part of an ABC transporter complex involved in host cell binding either through secretion of an adherence factor or through maintaining the architecture and integrity of the mycobacterial cell envelope (PubMed:17030567). Could be required for host endothelial-cell invasion and/or intracellular survival (PubMed:16586367). {ECO:0000269|PubMed:16586367, ECO:0000269|PubMed:17030567}.", 34),
    @("Rv1759c", 2, "wag22 wag22b Rv1759c MTCY28.25c", "", 34),
    @("Rv0849", 2, "Rv0849", "", 34),
    @("Rv1789", 2, "PPE26 Rv1789 LH57_09780", "FUNCTION: Probably plays a key role in regulating innate and adaptive immune responses through human Toll-like receptor 2 (TLR2). Interacts with TLR2, leading to the subsequent activation of the mitogen-activated protein kinase (MAPK) and nuclear factor kappa B (NF-kappa-B) signaling pathways. Stimulates macrophage activation by augmenting pro-inflammatory cytokine production (TNF-alpha, IL-6 and IL-12p40) and the expression of cell surface molecules (CD80, CD86, MHC class I and II). Also participates in adaptive immunity by directing Th1-polarised immune responses. {ECO:0000269|PubMed:26439698}.", 34),
    @("Rv1758", 2, "cut1 Rv1758", "FUNCTION: Does not exhibit cutinase activity. {ECO:0000269|PubMed:19225166}.", 34),
    @("Rv1761c", 1, "Rv1761c", "", 34),
    @("Rv0850", 1, "Rv0850", "", 34),
    @("Rv0228", 1, "Rv0228", "", 34),
    @("Rv1762c", 1, "Rv1762c", "", 34),
    @("Rv1081c", 1, "Rv1081c", "", 34),
    @("Rv2024c", 1, "Rv2024c", "", 34)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3]) {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
